$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 17:35"

# Row 4
$ws.Range("B4").Value = 1412661
$ws.Range("C4").Value = 4025
$ws.Range("D4").Value = 298979
$ws.Range("E4").Value = 1030018
$ws.Range("G4").Value = 239
$ws.Range("H4").Value = 83664

# Row 11
$ws.Range("B11").Value = 173647
$ws.Range("C11").Value = 476
$ws.Range("E11").Value = 17159
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = 7788

# Row 15
$ws.Range("B15").Value = 77729
$ws.Range("C15").Value = 3437
$ws.Range("D15").Value = 25977
$ws.Range("E15").Value = 49217
$ws.Range("G15").Value = 120
$ws.Range("H15").Value = 2535

# Row 17
$ws.Range("B17").Value = 71486
$ws.Range("C17").Value = 329
$ws.Range("E17").Value = 32235
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = 5209

# Row 20
$ws.Range("F20").Value = 432

# Row 52
$ws.Range("B52").Value = 8168
$ws.Range("C52").Value = 11
$ws.Range("E52").Value = 7907

# Row 57
$ws.Range("B57").Value = 6253
$ws.Range("C57").Value = 186
$ws.Range("D57").Value = 3058
$ws.Range("E57").Value = 2673
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 522

# Row 60
$ws.Range("D60").Value = 2408
$ws.Range("E60").Value = 2977

# Row 72
$ws.Range("A72").Value = "Grecia"
$ws.Range("B72").Value = 2760
$ws.Range("C72").Value = 16
$ws.Range("D72").Value = 1374
$ws.Range("E72").Value = 1231
$ws.Range("F72").Value = 28
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 155

# Row 73
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 2758
$ws.Range("C73").Value = 65
$ws.Range("D73").Value = 1789
$ws.Range("E73").Value = 934
$ws.Range("F73").Value = 30
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 35

# Row 83
$ws.Range("B83").Value = 1810
$ws.Range("C83").Value = 6
$ws.Range("D83").Value = 1326
$ws.Range("E83").Value = 405
$ws.Range("F83").Value = 7
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 79

# Row 93
$ws.Range("A93").Value = "Somalia"
$ws.Range("B93").Value = 1219
$ws.Range("C93").Value = 49
$ws.Range("D93").Value = 130
$ws.Range("E93").Value = 1037
$ws.Range("F93").Value = 2
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 52

# Row 94
$ws.Range("A94").Value = "Guatemala"
$ws.Range("B94").Value = 1199
$ws.Range("C94").Value = 85
$ws.Range("D94").Value = 120
$ws.Range("E94").Value = 1052
$ws.Range("F94").Value = 5
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 27

# Row 102
$ws.Range("A102").Value = "Republica de Chipre"
$ws.Range("B102").Value = 905
$ws.Range("C102").Value = 2
$ws.Range("D102").Value = 449
$ws.Range("E102").Value = 440
$ws.Range("F102").Value = 10
$ws.Range("H102").Value = 16

# Row 103
$ws.Range("A103").Value = "Maldivas"
$ws.Range("B103").Value = 904
$ws.Range("D103").Value = 29
$ws.Range("E103").Value = 872
$ws.Range("F103").Value = 2
$ws.Range("H103").Value = 3

# Row 121
$ws.Range("B121").Value = 582
$ws.Range("C121").Value = 6
$ws.Range("D121").Value = 392
$ws.Range("E121").Value = 181

# Row 128
$ws.Range("B128").Value = 439
$ws.Range("C128").Value = 2
$ws.Range("E128").Value = 85

# Row 131
$ws.Range("B131").Value = 372
$ws.Range("C131").Value = 15
$ws.Range("D131").Value = 78
$ws.Range("E131").Value = 252
$ws.Range("G131").Value = 2
$ws.Range("H131").Value = 42

# Row 134
$ws.Range("A134").Value = "Isla de Man"
$ws.Range("C134").Value = 1
$ws.Range("D134").Value = 272
$ws.Range("E134").Value = 37
$ws.Range("F134").Value = 20
$ws.Range("H134").Value = 23

# Row 135
$ws.Range("A135").Value = "Mauricio"
$ws.Range("B135").Value = 332
$ws.Range("D135").Value = 322
$ws.Range("E135").Value = 0
$ws.Range("F135").Value = 0
$ws.Range("H135").Value = 10

# Row 152
$ws.Range("B152").Value = 181
$ws.Range("C152").Value = 1
$ws.Range("D152").Value = 79
$ws.Range("E152").Value = 96
